$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.733350932598114
$ws.Range("B1").Value = 2.255666255950928
$ws.Range("D1").Value = 0.8924477100372314
$ws.Range("E1").Value = 0.8994224667549133
